# Edgar_scaling_mapping.xlsx update:
# - limit scaling for a few ISOs whose EDGAR data has jumps (rou, mkd, idn)
# - add select_scaling_year / start_scaling_year / end_scaling_year / Comment
#   columns to the "year" mapping sheet
# - update view state (active sheet, selections, frozen-pane position)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "year" sheet: new columns (E:H) and new rows (2:4)
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("year")

# Row 2 - rou / all, all years
$ws.Range("A2").Value = "rou"
$ws.Range("B2").Value = "all"
$ws.Range("C2").Value = "NA"
$ws.Range("D2").Value = "NA"
$ws.Range("E1").Value = "select_scaling_year"
$ws.Range("E2").Value = "1980, 1991,2000,2010"
$ws.Range("F1").Value = "start_scaling_year"
$ws.Range("G1").Value = "end_scaling_year"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "NA"

# Row 3 - mkd / all
$ws.Range("A3").Value = "mkd"
$ws.Range("B3").Value = "all"
$ws.Range("C3").Value = "NA"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "1980, 1992"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "NA"

# Row 4 - idn / 1A2, restricted start/end scaling year range
$ws.Range("A4").Value = "idn"
$ws.Range("B4").Value = "1A2"
$ws.Range("C4").Value = "NA"
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = 1990
$ws.Range("F4").Value = 2010
$ws.Range("G4").Value = "NA"

# Comment column
$ws.Range("H4").Value = "Eliminate jump in CO emissions present in EDGAR data"
$ws.Range("H1").Value = "Comment"
$ws.Range("H2").Value = "Reduce jumps in emissions"
$ws.Range("H3").Value = "Reduce jumps in emissions"

# widen the new "select_scaling_year" column
$ws.Columns("E:E").ColumnWidth = 17.67

# ------------------------------------------------------------------
# view state
# ------------------------------------------------------------------
# "map" sheet: scroll the frozen pane down and select D42 (no longer the
# active tab once "year" is selected below)
$wsMap = $wb.Worksheets.Item("map")
$wsMap.Range("D42").Select()

# "year" sheet becomes the active tab, with H3 selected
$ws.Activate()
$ws.Range("H3").Select()
